$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 9.861094666666666
$ws.Range("H2").Value = 29.583284
$ws.Range("I2").Value = 0.243709096397741
$ws.Range("J2").Value = 0.2437090963977409
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.04738633333333334
$ws.Range("N2").Value = 0.142159
$ws.Range("O2").Value = 0.05760194168856402
$ws.Range("P2").Value = 0.05760194168856402
$ws.Range("Q2").Value = 0.4672811189062223
$ws.Range("R2").Value = 4.205530070156
$ws.Range("S2").Value = 0.0140381171596753
$ws.Range("T2").Value = 0.0140381171596753

# Row 3
$ws.Range("G3").Value = 9.861094666666666
$ws.Range("H3").Value = 29.583284
$ws.Range("I3").Value = 0.243709096397741
$ws.Range("J3").Value = 0.2437090963977409
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6301496666666666
$ws.Range("N3").Value = 1.890449
$ws.Range("O3").Value = 0.7659981644722047
$ws.Range("P3").Value = 0.7659981644722047
$ws.Range("Q3").Value = 6.213965517168444
$ws.Range("R3").Value = 55.92568965451599
$ws.Range("S3").Value = 0.1866807205058492
$ws.Range("T3").Value = 0.1866807205058492

# Row 4
$ws.Range("G4").Value = 9.861094666666666
$ws.Range("H4").Value = 29.583284
$ws.Range("I4").Value = 0.243709096397741
$ws.Range("J4").Value = 0.2437090963977409
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1451156666666667
$ws.Range("N4").Value = 0.435347
$ws.Range("O4").Value = 0.1763998938392313
$ws.Range("P4").Value = 0.1763998938392313
$ws.Range("Q4").Value = 1.430999326616444
$ws.Range("R4").Value = 12.878993939548
$ws.Range("S4").Value = 0.04299025873221649
$ws.Range("T4").Value = 0.04299025873221649

# Row 5
$ws.Range("G5").Value = 29.35342966666667
$ws.Range("H5").Value = 88.060289
$ws.Range("I5").Value = 0.7254466225154019
$ws.Range("J5").Value = 0.7254466225154018
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.04738633333333334
$ws.Range("N5").Value = 0.142159
$ws.Range("O5").Value = 0.05760194168856402
$ws.Range("P5").Value = 0.05760194168856402
$ws.Range("Q5").Value = 1.390951402661222
$ws.Range("R5").Value = 12.518562623951
$ws.Range("S5").Value = 0.0417871340482979
$ws.Range("T5").Value = 0.04178713404829789

# Row 6
$ws.Range("G6").Value = 29.35342966666667
$ws.Range("H6").Value = 88.060289
$ws.Range("I6").Value = 0.7254466225154019
$ws.Range("J6").Value = 0.7254466225154018
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.6301496666666666
$ws.Range("N6").Value = 1.890449
$ws.Range("O6").Value = 0.7659981644722047
$ws.Range("P6").Value = 0.7659981644722047
$ws.Range("Q6").Value = 18.49705391997344
$ws.Range("R6").Value = 166.473485279761
$ws.Range("S6").Value = 0.5556907812693582
$ws.Range("T6").Value = 0.5556907812693581

# Row 7
$ws.Range("G7").Value = 29.35342966666667
$ws.Range("H7").Value = 88.060289
$ws.Range("I7").Value = 0.7254466225154019
$ws.Range("J7").Value = 0.7254466225154018
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.1451156666666667
$ws.Range("N7").Value = 0.435347
$ws.Range("O7").Value = 0.1763998938392313
$ws.Range("P7").Value = 0.1763998938392313
$ws.Range("Q7").Value = 4.259642515031445
$ws.Range("R7").Value = 38.336782635283
$ws.Range("S7").Value = 0.1279687071977458
$ws.Range("T7").Value = 0.1279687071977458

# Row 8
$ws.Range("G8").Value = 1.248038666666667
$ws.Range("H8").Value = 3.744116
$ws.Range("I8").Value = 0.03084428108685718
$ws.Range("J8").Value = 0.03084428108685716
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.04738633333333334
$ws.Range("N8").Value = 0.142159
$ws.Range("O8").Value = 0.05760194168856402
$ws.Range("P8").Value = 0.05760194168856402
$ws.Range("Q8").Value = 0.05913997627155556
$ws.Range("R8").Value = 0.5322597864440001
$ws.Range("S8").Value = 0.001776690480590825
$ws.Range("T8").Value = 0.001776690480590825

# Row 9
$ws.Range("G9").Value = 1.248038666666667
$ws.Range("H9").Value = 3.744116
$ws.Range("I9").Value = 0.03084428108685718
$ws.Range("J9").Value = 0.03084428108685716
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6301496666666666
$ws.Range("N9").Value = 1.890449
$ws.Range("O9").Value = 0.7659981644722047
$ws.Range("P9").Value = 0.7659981644722047
$ws.Range("Q9").Value = 0.7864511497871111
$ws.Range("R9").Value = 7.078060348084
$ws.Range("S9").Value = 0.02362666269699733
$ws.Range("T9").Value = 0.02362666269699733

# Row 10
$ws.Range("G10").Value = 1.248038666666667
$ws.Range("H10").Value = 3.744116
$ws.Range("I10").Value = 0.03084428108685718
$ws.Range("J10").Value = 0.03084428108685716
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.1451156666666667
$ws.Range("N10").Value = 0.435347
$ws.Range("O10").Value = 0.1763998938392313
$ws.Range("P10").Value = 0.1763998938392313
$ws.Range("Q10").Value = 0.1811099631391111
$ws.Range("R10").Value = 1.629989668252
$ws.Range("S10").Value = 0.005440927909269016
$ws.Range("T10").Value = 0.005440927909269013
